$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: mark as "waiting for payment" (highlighted Accent2 style) ---
$ws.Range("D15").Style = "Accent2"
$ws.Range("D15").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"
$ws.Range("E15").Style = "Accent2"

# --- Row 19: new receipt row (3D Hubs high speed prototype) ---
$ws.Range("A19").Value = "2018-04-10_3D_Hubs_High_Speed_Prototype.pdf"
$ws.Range("B19").Value = "High speed prototype 3D printing"
$ws.Range("C19").Value = "Thomas"
$ws.Range("D19").Value = 75.33

# --- Row 20: new receipt row (X-Axis motor/belt/pulley), waiting for payment ---
$ws.Range("B20").Value = "X-Axis high speed motor, more belt, more pulley"
$ws.Range("C20").Value = "Thomas"
$ws.Range("D20").Style = "Accent2"
$ws.Range("D20").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"
$ws.Range("E20").Style = "Accent2"
$ws.Range("E20").Value = "Waiting for payment to hit account for CAD value"

# --- Update selection to match the author's last-edited cell ---
$ws.Range("C20").Select()

Write-Output "done"
